# Auto-committed on 2022/03/11 週五
#
# 1) Insert a new row for table "TxArchiveTable" (歷史封存表設定檔) into the
#    "XX-系統" section, keeping the existing alphabetical ordering between
#    "TxApLogList" (row 328) and "TxAttachment" (previously row 329).
# 2) Refresh the "最後修改時間" (last modified) timestamps for two unrelated
#    tables: RepayActChangeLog (row 87) and AcDetail (row 135).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current TxAttachment row (329), pushing
# TxAttachment..TxUnLock (and everything below) down by one row, and
# inheriting the surrounding row's cell styles/formats.
$ws.Rows("329").Insert()

$ws.Range("A329").Value = "XX-系統"
$ws.Range("B329").Value = "TxArchiveTable"
$ws.Range("C329").Value = "歷史封存表設定檔"
$ws.Range("D329").Formula = '=HYPERLINK("[\\192.168.10.16\St1Share(NAS)\SKL\DB\GenTables\XX-系統\TxArchiveTable.xlsx]DBD!A1", "連結")'
$ws.Range("E329").Value = "2022年03月11日 12:05:24"

# Update last-modified timestamps for two other, unrelated tables.
$ws.Range("E87").Value = "2022年03月11日 11:23:55"
$ws.Range("E135").Value = "2022年03月11日 11:35:06"
